$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parts_1")

# Row 13 (KY-016 indicator LED) - supplier switched from Amazon to Banggood,
# with an updated price and an updated note about shipping from China.
$ws.Range("B13").Value = "Banggood"
$ws.Range("C13").Value = 4.27
$ws.Range("D13").Value = "ALLOW ~3 WEEKS TO ARRIVE FROM CHINA. Due to COVID, it's hard to find these in the US. You can get it quickly from Amazon B07KJYR8K1, but costs `$18."

# Update the selection to match what was left active in the saved file
$ws.Activate()
$ws.Range("C13").Select()
